# feat: add 2022-Q1 data
#
# The workbook has three sheets: "2021-Q1", "2021-Q4", "总计".
# This script:
#   1. Inserts a new "2022-Q1" sheet (a copy of the "2021-Q4" template)
#      right before the "总计" sheet, and fills it with the 2022-Q1
#      per-fund holdings data.
#   2. Inserts a new top data row into "总计" for the "2022-Q1" summary
#      line, shifting the existing rows down and renumbering the index
#      column.
#
# NOTE: worksheet references captured before an operation that changes
# sheet ordering/count (Copy, Add, Move, Delete) can become stale and
# point at the wrong sheet afterwards, so sheets are re-fetched by name
# right before they are used.

$wb = $excel.ActiveWorkbook

# --- 1. create the "2022-Q1" sheet from the "2021-Q4" template -------------
$ws2021Q4 = $wb.Worksheets.Item("2021-Q4")
$wsTotal  = $wb.Worksheets.Item("总计")
$ws2021Q4.Copy($wsTotal)

# sheet order shifted because of the Copy above - re-fetch before using
$wsTotal  = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($wsTotal.Index - 1)
$newSheet.Name = "2022-Q1"

# Cells B2:B3 and D2:G3 hold numeric-looking values that must stay text
# (matching the convention used on the other per-fund sheets), so force
# text formatting while writing them, then clear the formatting again so
# no stray style is left attached to the cells.
$newSheet  = $wb.Worksheets.Item("2022-Q1")
$rngCode = $newSheet.Range("B2:B3")
$rngNums = $newSheet.Range("D2:G3")
$rngCode.NumberFormat = "@"
$rngNums.NumberFormat = "@"

$newSheet.Range("B2").Value = "513030"
$newSheet.Range("C2").Value = "华安国际龙头(DAX)ETFQDII"
$newSheet.Range("D2").Value = "6.49"
$newSheet.Range("E2").Value = "92.80"
$newSheet.Range("F2").Value = "6.47"
$newSheet.Range("G2").Value = "0.4199"
$newSheet.Range("H2").Value = 4

$newSheet.Range("B3").Value = "006282"
$newSheet.Range("C3").Value = "上投摩根欧洲动力策略股票（QDII）"
$newSheet.Range("D3").Value = "0.48"
$newSheet.Range("E3").Value = "89.68"
$newSheet.Range("F3").Value = "2.01"
$newSheet.Range("G3").Value = "0.0096"
$newSheet.Range("H3").Value = 8

$rngCode.ClearFormats()
$rngNums.ClearFormats()

# --- 2. add the 2022-Q1 summary row to "总计" -------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Rows.Item(2).ClearFormats()

# Copy the index column's style onto the new A2 cell.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.43

# Renumber the index column for the rows that got pushed down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
